# Add two new columns (I: "I0", J: "IF") to the worksheet.
# I0 is a constant of 1 for every data row.
# IF mirrors the existing IP (column H) value for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting from the existing header cell (H1) onto the new
# header cells so they pick up the same style (bold, border, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows -----------------------------------------------------------
# IP values already present in column H (row -> value), re-used for the
# new IF column (J). I0 (column I) is always 1.
$ipValues = @{
    2  = 4
    3  = 5
    4  = 4
    5  = 5
    6  = 6
    7  = 3
    8  = 6
    9  = 5
    10 = 5
    11 = 6
    12 = 7
    13 = 5
    14 = 7
    15 = 6
    16 = 6
    17 = 6
    18 = 5
    19 = 3
    20 = 3
    21 = 6
    22 = 6
    23 = 6
    24 = 6
    25 = 6
    26 = 6
    27 = 7
    28 = 4
    29 = 7
    30 = 6
    31 = 3
    32 = 5
    33 = 3
}

foreach ($r in $ipValues.Keys) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValues[$r]
}
